$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>مربّاها</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>و</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>مشکلات</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>اقتصاد</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>بازه</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>باز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>در</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>برره</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>عدالت</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>قضا</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>یی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>برره</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>شرکت‌ها</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>هرم</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>در</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>برره</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>اعداد</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>هگزا</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>دس</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>مال</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>باز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ه</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>قل</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>دو</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>قل</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>در</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>برره</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>باز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>منطق</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>در</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>برره</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>نواح</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>سرو</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>س‌ده</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>سوپ</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>سبز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>جات</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>عوامل</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>اول</w:t></w:r></w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
